$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2..71: row, A (datetime serial), B (Lower2), C (MA)
# Exactly one of B/C is populated per row (alternating), matching the existing pattern.
$data = @(
    @(2, 44090.01041666666, 2.721714588703309, $null),
    @(3, 44090.66666666666, $null, 2.626344181182382),
    @(4, 44095.38541666666, 2.54601019334748, $null),
    @(5, 44095.6875, $null, 2.699081890485809),
    @(6, 44105.66666666666, 2.822815541294998, $null),
    @(7, 44105.9375, $null, 2.830988566391733),
    @(8, 44106.375, 2.576148539848103, $null),
    @(9, 44106.58333333334, $null, 2.630449014213091),
    @(10, 44108.09375, 2.460750845309257, $null),
    @(11, 44108.3125, $null, 2.548202028840201),
    @(12, 44110.45833333334, 2.445142707241562, $null),
    @(13, 44111.21875, $null, 2.221048679547392),
    @(14, 44112.29166666666, 2.039066820529677, $null),
    @(15, 44112.42708333334, $null, 2.182479925400636),
    @(16, 44120.16666666666, 2.144070585370677, $null),
    @(17, 44120.48958333334, $null, 2.191469830027393),
    @(18, 44123.15625, 2.110255916041447, $null),
    @(19, 44123.54166666666, $null, 2.107727062128928),
    @(20, 44124.32291666666, 1.972947284224925, $null),
    @(21, 44124.84375, $null, 1.915600803314232),
    @(22, 44129.875, 1.763241003717524, $null),
    @(23, 44130.08333333334, $null, 1.82061155675064),
    @(24, 44130.67708333334, 1.713633438908, $null),
    @(25, 44131.08333333334, $null, 1.708347335846987),
    @(26, 44132.29166666666, 1.523682772492656, $null),
    @(27, 44132.52083333334, $null, 1.599063398662677),
    @(28, 44137.375, 1.509468432597926, $null),
    @(29, 44137.90625, $null, 1.491771848317311),
    @(30, 44139.04166666666, 1.352453607204426, $null),
    @(31, 44139.41666666666, $null, 1.375994068739984),
    @(32, 44140.54166666666, 1.305564072633455, $null),
    @(33, 44140.58333333334, $null, 1.370658300155706),
    @(34, 44161.14583333334, 1.928037717057017, $null),
    @(35, 44161.85416666666, $null, 1.819467800485715),
    @(36, 44173.35416666666, 1.732894974855642, $null),
    @(37, 44173.67708333334, $null, 1.731193274007974),
    @(38, 44174.34375, 1.527500294467409, $null),
    @(39, 44174.48958333334, $null, 1.622281584551328),
    @(40, 44185.91666666666, 1.529514724800564, $null),
    @(41, 44186.07291666666, $null, 1.595449927742316),
    @(42, 44188.22916666666, 1.38851377236697, $null),
    @(43, 44189.0625, $null, 1.265586180957891),
    @(44, 44211.66666666666, 2.994398497615928, $null),
    @(45, 44211.9375, $null, 3.230715607712678),
    @(46, 44217.9375, 2.960696161077601, $null),
    @(47, 44218.15625, $null, 3.101020120276233),
    @(48, 44233.07291666666, 6.117283855637567, $null),
    @(49, 44233.73958333334, $null, 6.022875422410603),
    @(50, 44237.53125, 7.370374437788142, $null),
    @(51, 44237.76041666666, $null, 8.069859602960564),
    @(52, 44242.09375, 8.037609503087824, $null),
    @(53, 44242.35416666666, $null, 8.369873372139002),
    @(54, 44249.59375, 8.961131125366157, $null),
    @(55, 44249.61458333334, $null, 10.28836147115143),
    @(56, 44267.83333333334, 14.12585757239154, $null),
    @(57, 44268.21875, $null, 14.49075860085157),
    @(58, 44279.92708333334, 13.73668846338735, $null),
    @(59, 44280.77083333334, $null, 13.08219313372434),
    @(60, 44304.05208333334, 23.54699873187199, $null),
    @(61, 44304.25, $null, 24.35933694351101),
    @(62, 44326.83333333334, 40.39098955608101, $null),
    @(63, 44326.90625, $null, 43.53030316279118),
    @(64, 44335.52083333334, 36.46550974054149, $null),
    @(65, 44335.63541666666, $null, 45.68648690290231),
    @(66, 44339.40625, 21.93368251561866, $null),
    @(67, 44339.83333333334, $null, 23.61274779655797),
    @(68, 44367.5, 31.64069288468217, $null),
    @(69, 44367.63541666666, $null, 33.99281598173719),
    @(70, 44432.67708333334, 68.63250106733484, $null),
    @(71, 44433.11458333334, $null, 71.89387507720097)
)

# Extend formatting (date style incl. number format/border/font/alignment) from A2
# down through the new rows (22..71) before writing values, so the new date cells
# match the existing column-A style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A22:A71").PasteSpecial(-4122) | Out-Null

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $c = $row[3]

    $ws.Cells.Item($r, 1).Value = $a

    if ($null -ne $b) {
        $ws.Cells.Item($r, 2).Value = $b
    }
    if ($null -ne $c) {
        $ws.Cells.Item($r, 3).Value = $c
    }
}
